$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 43, shifting existing rows 43:93 down to 44:94
$ws.Rows.Item(43).Insert()

# Populate the new row 43 with values.
# Columns A, B, C, E, F, G, H, I, J, L, M keep the same values the old row 43 had.
$ws.Range("A43").Value = 11
$ws.Range("B43").Value = "Vega Monumental Concepción"
$ws.Range("C43").Value = "Bíobío"
$ws.Range("D43").Value = 44554
$ws.Range("E43").Value = 8
$ws.Range("F43").Value = "Fruta"
$ws.Range("G43").Value = 100109
$ws.Range("H43").Value = "Uva"
$ws.Range("I43").Value = 100109001
$ws.Range("J43").Value = "Uva"
$ws.Range("K43").Value = "Superior Seedless"
$ws.Range("L43").Value = "Primera"
$ws.Range("M43").Value = 100
$ws.Range("N43").Value = 13000
$ws.Range("O43").Value = 14000
$ws.Range("P43").Value = 13500
$ws.Range("Q43").Value = "`$/bandeja 10 kilos"
$ws.Range("R43").Value = "Provincia del Elquí"
$ws.Range("S43").Value = 1350
$ws.Range("T43").Value = 10

# Apply the same date number format used on other cells in column D
$ws.Range("D43").NumberFormat = "YYYY-MM-DD HH:MM:SS"
